$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr_B = New-Object 'object[,]' 24,4
$arr_B[0,0] = 11.21885214680199
$arr_B[0,1] = 8.686567973578404
$arr_B[0,2] = 14.7667239225993
$arr_B[0,3] = 15.86351786790864
$arr_B[1,0] = 11.08944375622857
$arr_B[1,1] = 8.589577762117903
$arr_B[1,2] = 14.74302752984397
$arr_B[1,3] = 15.86070817958923
$arr_B[2,0] = 11.01282307159478
$arr_B[2,1] = 8.532151109116864
$arr_B[2,2] = 14.73140451495517
$arr_B[2,3] = 15.86175528543909
$arr_B[3,0] = 10.9823534881586
$arr_B[3,1] = 8.509311558298387
$arr_B[3,2] = 14.72740706274803
$arr_B[3,3] = 15.8628799333422
$arr_B[4,0] = 10.97734073203774
$arr_B[4,1] = 8.505553814540717
$arr_B[4,2] = 14.72678800310026
$arr_B[4,3] = 15.86310885275589
$arr_B[5,0] = 11.01240904280469
$arr_B[5,1] = 8.531840774593849
$arr_B[5,2] = 14.73134760798726
$arr_B[5,3] = 15.86176762601152
$arr_B[6,0] = 11.17366828072594
$arr_B[6,1] = 8.652701568022952
$arr_B[6,2] = 14.75794799145397
$arr_B[6,3] = 15.86197429030758
$arr_B[7,0] = 11.5105460084744
$arr_B[7,1] = 8.905301536562151
$arr_B[7,2] = 14.83317459888256
$arr_B[7,3] = 15.8843288619208
$arr_B[8,0] = 11.76807734509981
$arr_B[8,1] = 9.098625853743021
$arr_B[8,2] = 14.90225828687454
$arr_B[8,3] = 15.91405428080753
$arr_B[9,0] = 11.88685799676932
$arr_B[9,1] = 9.1878689700085
$arr_B[9,2] = 14.93662224776484
$arr_B[9,3] = 15.93044202220468
$arr_B[10,0] = 11.93202502610024
$arr_B[10,1] = 9.221817365438358
$arr_B[10,2] = 14.95005107366441
$arr_B[10,3] = 15.93705709762867
$arr_B[11,0] = 11.92228990373183
$arr_B[11,1] = 9.214499643534955
$arr_B[11,2] = 14.9471405452031
$arr_B[11,3] = 15.93561426081444
$arr_B[12,0] = 11.89057042624893
$arr_B[12,1] = 9.190659026437343
$arr_B[12,2] = 14.93771874818942
$arr_B[12,3] = 15.93097806001422
$arr_B[13,0] = 11.87116435259194
$arr_B[13,1] = 9.176075040447712
$arr_B[13,2] = 14.93200159230131
$arr_B[13,3] = 15.92819148359225
$arr_B[14,0] = 11.76034350695143
$arr_B[14,1] = 9.092816934085114
$arr_B[14,2] = 14.90007108994624
$arr_B[14,3] = 15.91304074459513
$arr_B[15,0] = 11.69274214375854
$arr_B[15,1] = 9.042050065302364
$arr_B[15,2] = 14.88123068367749
$arr_B[15,3] = 15.90447839193422
$arr_B[16,0] = 11.65401567609672
$arr_B[16,1] = 9.012974525092348
$arr_B[16,2] = 14.87067108761232
$arr_B[16,3] = 15.89982344800367
$arr_B[17,0] = 11.64093179396634
$arr_B[17,1] = 9.00315238443998
$arr_B[17,2] = 14.86714353704467
$arr_B[17,3] = 15.89829379432351
$arr_B[18,0] = 11.69992261426439
$arr_B[18,1] = 9.047441667187408
$arr_B[18,2] = 14.88320766334133
$arr_B[18,3] = 15.90536195411543
$arr_B[19,0] = 11.89988249055835
$arr_B[19,1] = 9.197657675466521
$arr_B[19,2] = 14.94047492399
$arr_B[19,3] = 15.93232873661383
$arr_B[20,0] = 12.03163873195546
$arr_B[20,1] = 9.296714923758053
$arr_B[20,2] = 14.98032343846358
$arr_B[20,3] = 15.95233787749979
$arr_B[21,0] = 11.96123526122855
$arr_B[21,1] = 9.243776239650382
$arr_B[21,2] = 14.9588362803042
$arr_B[21,3] = 15.9414413790735
$arr_B[22,0] = 11.69667588867921
$arr_B[22,1] = 9.045003775680945
$arr_B[22,2] = 14.88231302309567
$arr_B[22,3] = 15.90496166170436
$arr_B[23,0] = 11.41746291798102
$arr_B[23,1] = 8.835477300056414
$arr_B[23,2] = 14.8103778025958
$arr_B[23,3] = 15.87593843228408
$ws.Range("B2:E25").Value2 = $arr_B

$arr_G = New-Object 'object[,]' 24,5
$arr_G[0,0] = 58.36701257356254
$arr_G[0,1] = 21.52826387766281
$arr_G[0,2] = 32.09379621313019
$arr_G[0,3] = 9.272602774495139
$arr_G[0,4] = 12.56488108653791
$arr_G[1,0] = 58.11686769038016
$arr_G[1,1] = 21.52605082508355
$arr_G[1,2] = 32.06681038088814
$arr_G[1,3] = 9.288739119551714
$arr_G[1,4] = 12.49240571330491
$arr_G[2,0] = 57.97366579933315
$arr_G[2,1] = 21.52761520940386
$arr_G[2,2] = 32.05471715950679
$arr_G[2,3] = 9.299493341223537
$arr_G[2,4] = 12.45116407187271
$arr_G[3,0] = 57.91795368835001
$arr_G[3,1] = 21.5289880277403
$arr_G[3,2] = 32.05091602523354
$arr_G[3,3] = 9.304088931241854
$arr_G[3,4] = 12.43519324518289
$arr_G[4,0] = 57.9088632728887
$arr_G[4,1] = 21.52926039718046
$arr_G[4,2] = 32.05035294429135
$arr_G[4,3] = 9.304864908821211
$arr_G[4,4] = 12.43259220593533
$arr_G[5,0] = 57.9729037024578
$arr_G[5,1] = 21.52763074634323
$arr_G[5,2] = 32.05466133176841
$arr_G[5,3] = 9.299554455479749
$arr_G[5,4] = 12.45094528067074
$arr_G[6,0] = 58.2786265666921
$arr_G[6,1] = 21.52689432407542
$arr_G[6,2] = 32.08356344178465
$arr_G[6,3] = 9.277991095825717
$arr_G[6,4] = 12.53922533312701
$arr_G[7,0] = 58.95903612835755
$arr_G[7,1] = 21.54862649306337
$arr_G[7,2] = 32.17567880608701
$arr_G[7,3] = 9.242407739034284
$arr_G[7,4] = 12.73742732394209
$arr_G[8,0] = 59.50598403543066
$arr_G[8,1] = 21.57868760947462
$arr_G[8,2] = 32.26482959508747
$arr_G[8,3] = 9.220331935495272
$arr_G[8,4] = 12.89722441947086
$arr_G[9,0] = 59.76447746000818
$arr_G[9,1] = 21.59541117084907
$arr_G[9,2] = 32.310010150624
$arr_G[9,3] = 9.211168309082646
$arr_G[9,4] = 12.97273967259518
$arr_G[10,0] = 59.86370432606627
$arr_G[10,1] = 21.60218071374457
$arr_G[10,2] = 32.32777951293585
$arr_G[10,3] = 9.207824334477714
$arr_G[10,4] = 13.00171817150729
$arr_G[11,0] = 59.84227524982953
$arr_G[11,1] = 21.60070337942437
$arr_G[11,2] = 32.32392326883531
$arr_G[11,3] = 9.208538915637895
$arr_G[11,4] = 12.9954605006598
$arr_G[12,0] = 59.77261434204144
$arr_G[12,1] = 21.59595936392841
$arr_G[12,2] = 32.31145883881024
$arr_G[12,3] = 9.210890672781796
$arr_G[12,4] = 12.97511622146351
$arr_G[13,0] = 59.73011809246938
$arr_G[13,1] = 21.59311033277797
$arr_G[13,2] = 32.30390987938536
$arr_G[13,3] = 9.212347604680764
$arr_G[13,4] = 12.96270387125709
$arr_G[14,0] = 59.48928140204438
$arr_G[14,1] = 21.57765589285172
$arr_G[14,2] = 32.26196962950193
$arr_G[14,3] = 9.220948459682196
$arr_G[14,4] = 12.89234413711208
$arr_G[15,0] = 59.3439803530733
$arr_G[15,1] = 21.56895480046322
$arr_G[15,2] = 32.23742219326751
$arr_G[15,3] = 9.226449684380917
$arr_G[15,4] = 12.84988727558931
$arr_G[16,0] = 59.26132179580005
$arr_G[16,1] = 21.56423723056723
$arr_G[16,2] = 32.22373883924698
$arr_G[16,3] = 9.229696573459815
$arr_G[16,4] = 12.82573435226155
$arr_G[17,0] = 59.23349366917593
$arr_G[17,1] = 21.56268929379875
$arr_G[17,2] = 32.21918085204362
$arr_G[17,3] = 9.230810131632355
$arr_G[17,4] = 12.81760317030674
$arr_G[18,0] = 59.3593535885355
$arr_G[18,1] = 21.56985134743124
$arr_G[18,2] = 32.23999025014788
$arr_G[18,3] = 9.22585550934056
$arr_G[18,4] = 12.85437940019163
$arr_G[19,0] = 59.79303946746447
$arr_G[19,1] = 21.59734095942191
$arr_G[19,2] = 32.31510205923001
$arr_G[19,3] = 9.210196484789931
$arr_G[19,4] = 12.98108163798219
$arr_G[20,0] = 60.08426913050467
$arr_G[20,1] = 21.61785140978197
$arr_G[20,2] = 32.36803909921986
$arr_G[20,3] = 9.200697239297607
$arr_G[20,4] = 13.06610714364732
$arr_G[21,0] = 59.9281391671209
$arr_G[21,1] = 21.60667241511943
$arr_G[21,2] = 32.33943528296598
$arr_G[21,3] = 9.205700018818114
$arr_G[21,4] = 13.02053235710536
$arr_G[22,0] = 59.35240061306639
$arr_G[22,1] = 21.56944513117596
$arr_G[22,2] = 32.23882789393967
$arr_G[22,3] = 9.226123873692888
$arr_G[22,4] = 12.85234771115952
$arr_G[23,0] = 58.76652866119133
$arr_G[23,1] = 21.54026867046421
$arr_G[23,2] = 32.14697640620047
$arr_G[23,3] = 9.251318374247074
$arr_G[23,4] = 12.68123357233142
$ws.Range("G2:K25").Value2 = $arr_G

$arr_M = New-Object 'object[,]' 24,1
$arr_M[0,0] = 18.23478337878501
$arr_M[1,0] = 18.21854061666204
$arr_M[2,0] = 18.21230440263233
$arr_M[3,0] = 18.21070574239086
$arr_M[4,0] = 18.21049729482783
$arr_M[5,0] = 18.21227902232895
$arr_M[6,0] = 18.22840923473713
$arr_M[7,0] = 18.28953676508974
$arr_M[8,0] = 18.35218696394463
$arr_M[9,0] = 18.38447620876891
$arr_M[10,0] = 18.39724177353927
$arr_M[11,0] = 18.39446864945039
$arr_M[12,0] = 18.38551569032996
$arr_M[13,0] = 18.38010164534134
$arr_M[14,0] = 18.35015248880963
$arr_M[15,0] = 18.33274570667766
$arr_M[16,0] = 18.32309088976094
$arr_M[17,0] = 18.31988345818839
$arr_M[18,0] = 18.33456177518849
$arr_M[19,0] = 18.38813083670576
$arr_M[20,0] = 18.4262756413708
$arr_M[21,0] = 18.40563257473384
$arr_M[22,0] = 18.33373963186823
$arr_M[23,0] = 18.2698649116665
$ws.Range("M2:M25").Value2 = $arr_M
